$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.763.35"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "1.675.99"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9965"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.82"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4635"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2591"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06144"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "1.669.04"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06979"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.95"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5820"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.358"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.42"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9981"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9982"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "25.752.95"
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006717"
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "1.880.29"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.453"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.659"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.247"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.62"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.01"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.385"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.730"
$ws.Range("E28").Value = "  +5.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.65"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.953"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07674"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.612"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04367"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.600"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6125"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9514"
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9221"
$ws.Range("E37").Value = "  +7.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.465"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "108.64"
$ws.Range("E39").Value = "  +10.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9974"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.863"
$ws.Range("E41").Value = "  +5.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01455"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.087"
$ws.Range("E43").Value = "  +9.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3731"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1117"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05285"
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.149"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "31.10"
$ws.Range("E48").Value = "  +8.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.637"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.208"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9982"
$ws.Range("E51").Value = "  -0.36%  "
